$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.929.05"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.358.44"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "504.79"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "130.49"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("D9").Value = "2.375.29"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "4.81"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "0.322"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "2.779.22"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "55.810.78"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "21.55"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.366.15"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "9.92"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "310.54"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "4.02"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "65.35"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "7.15"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "170.64"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "0.0₃0707"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "36.15"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("D44").Value = "126.30"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("D45").Value = "0.559"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "0.0898"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "240.89"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "16.90"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "16.71"
$ws.Range("E51").Value = "  -2.45%  "
